# ---------------------------------------------------------------------------
# Update "latest output" optimisation results (run 265)
#
# 1) Schedule sheet: revise the pumping-window cost/unit-cost figures and
#    append a second scheduled pumping block (row 3).
# 2) Detailed sheet: re-classify several half-hours from forecast to
#    historical (price updates as the actuals land) and append a full new
#    day (rows 50-97) of forecast price/pump-status data.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Schedule sheet -------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E2").Value = 129.72843975
$schedule.Range("F2").Value = 3.119972095959596

# New pumping block appended as row 3 (copy the datetime format used by
# row 2 for the Start Time / Stop Time columns).
$schedule.Range("A3").Value = 46082.22916666666
$schedule.Range("A3").NumberFormat = $schedule.Range("A2").NumberFormat()
$schedule.Range("B3").Value = 46082.6875
$schedule.Range("B3").NumberFormat = $schedule.Range("B2").NumberFormat()
$schedule.Range("C3").Value = 11
$schedule.Range("D3").Value = 41.58
$schedule.Range("E3").Value = 240.993558
$schedule.Range("F3").Value = 5.795900865800866

# --- Detailed sheet ---------------------------------------------------------
$detailed = $wb.Worksheets.Item("Detailed")

# Price corrections / forecast -> historical re-classification for the
# half-hours that have now occurred (rows 13-47).
$detailed.Range("B13").Value = 65
$detailed.Range("B14").Value = 57.36
$detailed.Range("C15").Value = "historical"
$detailed.Range("B16").Value = 56.98; $detailed.Range("C16").Value = "historical"
$detailed.Range("C17").Value = "historical"
$detailed.Range("B18").Value = 1.16971; $detailed.Range("C18").Value = "historical"
$detailed.Range("B19").Value = 1.15893; $detailed.Range("C19").Value = "historical"
$detailed.Range("B20").Value = 1.16138; $detailed.Range("C20").Value = "historical"
$detailed.Range("B21").Value = 1.13838; $detailed.Range("C21").Value = "historical"
$detailed.Range("B22").Value = 1.07591; $detailed.Range("C22").Value = "historical"
$detailed.Range("B23").Value = 0.7; $detailed.Range("C23").Value = "historical"
$detailed.Range("B24").Value = 0.51; $detailed.Range("C24").Value = "historical"
$detailed.Range("B25").Value = 0.7; $detailed.Range("C25").Value = "historical"
$detailed.Range("C26").Value = "historical"
$detailed.Range("B27").Value = 0.51; $detailed.Range("C27").Value = "historical"
$detailed.Range("B28").Value = 0.7; $detailed.Range("C28").Value = "historical"
$detailed.Range("B29").Value = 0.7; $detailed.Range("C29").Value = "historical"
$detailed.Range("C30").Value = "historical"
$detailed.Range("B31").Value = -10.08057; $detailed.Range("C31").Value = "historical"
$detailed.Range("B32").Value = 0.5099399999999999; $detailed.Range("C32").Value = "historical"
$detailed.Range("B33").Value = 0.51
$detailed.Range("B34").Value = -3.75989
$detailed.Range("B35").Value = -5.01
$detailed.Range("B36").Value = 27.68107
$detailed.Range("B38").Value = 38.7
$detailed.Range("B39").Value = 46.78124
$detailed.Range("B40").Value = 56.65075
$detailed.Range("B41").Value = 57.36
$detailed.Range("B42").Value = 57.32
$detailed.Range("B43").Value = 54.67362
$detailed.Range("B44").Value = 54.26951
$detailed.Range("B45").Value = 56.98
$detailed.Range("B46").Value = 56.98
$detailed.Range("B47").Value = 57.06

# Append a full new day of forecast price / pump-status data (rows 50-97).
$newRows = @(
    @(50, 46082, 56.98, "forecast", 46082, "OFF"),
    @(51, 46082.02083333334, 57.06, "forecast", 46082, "OFF"),
    @(52, 46082.04166666666, 54.41385, "forecast", 46082, "OFF"),
    @(53, 46082.0625, 40.22357, "forecast", 46082, "OFF"),
    @(54, 46082.08333333334, 37.89, "forecast", 46082, "OFF"),
    @(55, 46082.10416666666, 37.89, "forecast", 46082, "OFF"),
    @(56, 46082.125, 37.89, "forecast", 46082, "OFF"),
    @(57, 46082.14583333334, 35.88, "forecast", 46082, "OFF"),
    @(58, 46082.16666666666, 35.88, "forecast", 46082, "OFF"),
    @(59, 46082.1875, 37.89, "forecast", 46082, "OFF"),
    @(60, 46082.20833333334, 37.89, "forecast", 46082, "OFF"),
    @(61, 46082.22916666666, 37.89, "forecast", 46082, "ON"),
    @(62, 46082.25, 37.89, "forecast", 46082, "ON"),
    @(63, 46082.27083333334, 37.89, "forecast", 46082, "ON"),
    @(64, 46082.29166666666, 37.89, "forecast", 46082, "ON"),
    @(65, 46082.3125, 0.51, "forecast", 46082, "ON"),
    @(66, 46082.33333333334, 0.05359, "forecast", 46082, "ON"),
    @(67, 46082.35416666666, 0.51, "forecast", 46082, "ON"),
    @(68, 46082.375, 0.5099399999999999, "forecast", 46082, "ON"),
    @(69, 46082.39583333334, 0.51, "forecast", 46082, "ON"),
    @(70, 46082.41666666666, 0.01132, "forecast", 46082, "ON"),
    @(71, 46082.4375, -0.10591, "forecast", 46082, "ON"),
    @(72, 46082.45833333334, -0.10246, "forecast", 46082, "ON"),
    @(73, 46082.47916666666, -0.10158, "forecast", 46082, "ON"),
    @(74, 46082.5, 0.01016, "forecast", 46082, "ON"),
    @(75, 46082.52083333334, 1.70773, "forecast", 46082, "ON"),
    @(76, 46082.54166666666, 1.73693, "forecast", 46082, "ON"),
    @(77, 46082.5625, 0.50996, "forecast", 46082, "ON"),
    @(78, 46082.58333333334, 0.50996, "forecast", 46082, "ON"),
    @(79, 46082.60416666666, 1.74228, "forecast", 46082, "ON"),
    @(80, 46082.625, 13.83096, "forecast", 46082, "ON"),
    @(81, 46082.64583333334, 35.88, "forecast", 46082, "ON"),
    @(82, 46082.66666666666, 37.89, "forecast", 46082, "ON"),
    @(83, 46082.6875, 40.64438, "forecast", 46082, "OFF"),
    @(84, 46082.70833333334, 50.96614, "forecast", 46082, "OFF"),
    @(85, 46082.72916666666, 52.6633, "forecast", 46082, "OFF"),
    @(86, 46082.75, 57.31, "forecast", 46082, "OFF"),
    @(87, 46082.77083333334, 57.36, "forecast", 46082, "OFF"),
    @(88, 46082.79166666666, 57.36, "forecast", 46082, "OFF"),
    @(89, 46082.8125, 59.19004, "forecast", 46082, "OFF"),
    @(90, 46082.83333333334, 57.31, "forecast", 46082, "OFF"),
    @(91, 46082.85416666666, 57.31, "forecast", 46082, "OFF"),
    @(92, 46082.875, 57.0601, "forecast", 46082, "OFF"),
    @(93, 46082.89583333334, 57.06, "forecast", 46082, "OFF"),
    @(94, 46082.91666666666, 53.23234, "forecast", 46082, "OFF"),
    @(95, 46082.9375, 57.06, "forecast", 46082, "OFF"),
    @(96, 46082.95833333334, 52.57899, "forecast", 46082, "OFF"),
    @(97, 46082.97916666666, 52.45471, "forecast", 46082, "OFF")
)

foreach ($r in $newRows) {
    $row = $r[0]
    $detailed.Range("A$row").Value = $r[1]
    $detailed.Range("A$row").NumberFormat = $detailed.Range("A49").NumberFormat()
    $detailed.Range("B$row").Value = $r[2]
    $detailed.Range("C$row").Value = $r[3]
    $detailed.Range("D$row").Value = $r[4]
    $detailed.Range("D$row").NumberFormat = $detailed.Range("D49").NumberFormat()
    $detailed.Range("E$row").Value = $r[5]
}
